$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

$tbl.Cell(1, 1).Range.Text = "73×11=803"
$tbl.Cell(1, 2).Range.Text = "55×60=3300"
$tbl.Cell(1, 3).Range.Text = "58×61=3538"
$tbl.Cell(1, 4).Range.Text = "84×49=4116"
$tbl.Cell(1, 5).Range.Text = "94×39=3666"
$tbl.Cell(5, 1).Range.Text = "95×33=3135"
$tbl.Cell(5, 2).Range.Text = "41×26=1066"
$tbl.Cell(5, 3).Range.Text = "80×39=3120"
$tbl.Cell(5, 4).Range.Text = "40×23=920"
$tbl.Cell(5, 5).Range.Text = "83×65=5395"
$tbl.Cell(10, 1).Range.Text = "15×33=495"
$tbl.Cell(10, 2).Range.Text = "71×74=5254"
$tbl.Cell(10, 3).Range.Text = "66×84=5544"
$tbl.Cell(10, 4).Range.Text = "53×79=4187"
$tbl.Cell(10, 5).Range.Text = "28×14=392"
$tbl.Cell(15, 1).Range.Text = "86×64=5504"
$tbl.Cell(15, 2).Range.Text = "59×21=1239"
$tbl.Cell(15, 3).Range.Text = "77×99=7623"
$tbl.Cell(15, 4).Range.Text = "79×15=1185"
$tbl.Cell(15, 5).Range.Text = "41×61=2501"
$tbl.Cell(20, 1).Range.Text = "15×97=1455"
$tbl.Cell(20, 2).Range.Text = "40×92=3680"
$tbl.Cell(20, 3).Range.Text = "53×17=901"
$tbl.Cell(20, 4).Range.Text = "29×70=2030"
$tbl.Cell(20, 5).Range.Text = "87×88=7656"

Write-Output "Updated 25 cells"
